$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Bdnf"
$ws.Cells.Item(2, 3).Value = "Ngfr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.03885866666666667
$ws.Cells.Item(2, 8).Value = 0.116576
$ws.Cells.Item(2, 9).Value = 0.01924839521029073
$ws.Cells.Item(2, 10).Value = 0.01924839521029073
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.261293
$ws.Cells.Item(2, 14).Value = 0.783879
$ws.Cells.Item(2, 15).Value = 0.04010308082688332
$ws.Cells.Item(2, 16).Value = 0.04010308082688332
$ws.Cells.Item(2, 17).Value = 0.01015349758933333
$ws.Cells.Item(2, 18).Value = 0.09138147830399999
$ws.Cells.Item(2, 19).Value = 0.0007719199489060831
$ws.Cells.Item(2, 20).Value = 0.0007719199489060829
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Bdnf"
$ws.Cells.Item(3, 3).Value = "Ngfr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.03885866666666667
$ws.Cells.Item(3, 8).Value = 0.116576
$ws.Cells.Item(3, 9).Value = 0.01924839521029073
$ws.Cells.Item(3, 10).Value = 0.01924839521029073
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.311698333333333
$ws.Cells.Item(3, 14).Value = 3.935095
$ws.Cells.Item(3, 15).Value = 0.201318612753326
$ws.Cells.Item(3, 16).Value = 0.201318612753326
$ws.Cells.Item(3, 17).Value = 0.05097084830222223
$ws.Cells.Item(3, 18).Value = 0.45873763472
$ws.Cells.Item(3, 19).Value = 0.003875060221463494
$ws.Cells.Item(3, 20).Value = 0.003875060221463494
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Bdnf"
$ws.Cells.Item(4, 3).Value = "Ngfr"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.03885866666666667
$ws.Cells.Item(4, 8).Value = 0.116576
$ws.Cells.Item(4, 9).Value = 0.01924839521029073
$ws.Cells.Item(4, 10).Value = 0.01924839521029073
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.942543
$ws.Cells.Item(4, 14).Value = 14.827629
$ws.Cells.Item(4, 15).Value = 0.7585783064197906
$ws.Cells.Item(4, 16).Value = 0.7585783064197906
$ws.Cells.Item(4, 17).Value = 0.1920606309226666
$ws.Cells.Item(4, 18).Value = 1.728545678304
$ws.Cells.Item(4, 19).Value = 0.01460141503992115
$ws.Cells.Item(4, 20).Value = 0.01460141503992115
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Bdnf"
$ws.Cells.Item(5, 3).Value = "Ngfr"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.979941666666667
$ws.Cells.Item(5, 8).Value = 5.939825
$ws.Cells.Item(5, 9).Value = 0.9807516047897092
$ws.Cells.Item(5, 10).Value = 0.9807516047897092
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.261293
$ws.Cells.Item(5, 14).Value = 0.783879
$ws.Cells.Item(5, 15).Value = 0.04010308082688332
$ws.Cells.Item(5, 16).Value = 0.04010308082688332
$ws.Cells.Item(5, 17).Value = 0.5173448979083334
$ws.Cells.Item(5, 18).Value = 4.656104081175
$ws.Cells.Item(5, 19).Value = 0.03933116087797724
$ws.Cells.Item(5, 20).Value = 0.03933116087797724
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Bdnf"
$ws.Cells.Item(6, 3).Value = "Ngfr"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.979941666666667
$ws.Cells.Item(6, 8).Value = 5.939825
$ws.Cells.Item(6, 9).Value = 0.9807516047897092
$ws.Cells.Item(6, 10).Value = 0.9807516047897092
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.311698333333333
$ws.Cells.Item(6, 14).Value = 3.935095
$ws.Cells.Item(6, 15).Value = 0.201318612753326
$ws.Cells.Item(6, 16).Value = 0.201318612753326
$ws.Cells.Item(6, 17).Value = 2.597086184263889
$ws.Cells.Item(6, 18).Value = 23.373775658375
$ws.Cells.Item(6, 19).Value = 0.1974435525318625
$ws.Cells.Item(6, 20).Value = 0.1974435525318625
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Bdnf"
$ws.Cells.Item(7, 3).Value = "Ngfr"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.979941666666667
$ws.Cells.Item(7, 8).Value = 5.939825
$ws.Cells.Item(7, 9).Value = 0.9807516047897092
$ws.Cells.Item(7, 10).Value = 0.9807516047897092
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.942543
$ws.Cells.Item(7, 14).Value = 14.827629
$ws.Cells.Item(7, 15).Value = 0.7585783064197906
$ws.Cells.Item(7, 16).Value = 0.7585783064197906
$ws.Cells.Item(7, 17).Value = 9.785946824991667
$ws.Cells.Item(7, 18).Value = 88.073521424925
$ws.Cells.Item(7, 19).Value = 0.7439768913798694
$ws.Cells.Item(7, 20).Value = 0.7439768913798694
